$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(8).Cut($ws.Rows.Item(4))
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
Write-Output "--- final ---"
for ($r=1; $r -le 8; $r++) {
  $a = $ws.Cells.Item($r,1).Value()
  $b = $ws.Cells.Item($r,2).Value()
  Write-Output "Row $r : A=$a | B=$b"
}
